# Add data for 2024-09-02
# Updates the 2024 (K column) cumulative totals across the Citywide Totals,
# By Neighborhood, and individual neighborhood sheets to reflect the newly
# added day of violent crime data.

$wb = $excel.ActiveWorkbook


# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5379
$ws.Range("K3").Value = 5537
$ws.Range("K4").Value = 1156
$ws.Range("K5").Value = 397
$ws.Range("K6").Value = 6165
$ws.Range("K7").Value = 18634

# Norwood Park
$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 41

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 237

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 372
$ws.Range("K4").Value = 71
$ws.Range("K6").Value = 420
$ws.Range("K7").Value = 1241

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 141
$ws.Range("K3").Value = 150
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 416

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 294
$ws.Range("K6").Value = 236
$ws.Range("K7").Value = 801

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 105
$ws.Range("K7").Value = 315

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K6").Value = 183
$ws.Range("K7").Value = 627

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 107
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 423

# Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 314

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K5").Value = 44
$ws.Range("K7").Value = 554
$ws.Range("K8").Value = 1241
$ws.Range("K9").Value = 75
$ws.Range("K10").Value = 104
$ws.Range("K11").Value = 353
$ws.Range("K14").Value = 98
$ws.Range("K18").Value = 126
$ws.Range("K19").Value = 547
$ws.Range("K23").Value = 194
$ws.Range("K29").Value = 1001
$ws.Range("K33").Value = 801
$ws.Range("K35").Value = 29
$ws.Range("K37").Value = 627
$ws.Range("K39").Value = 23
$ws.Range("K42").Value = 690
$ws.Range("K43").Value = 165
$ws.Range("K44").Value = 157
$ws.Range("K47").Value = 127
$ws.Range("K48").Value = 235
$ws.Range("K52").Value = 488
$ws.Range("K53").Value = 237
$ws.Range("K55").Value = 208
$ws.Range("K61").Value = 16
$ws.Range("K63").Value = 54
$ws.Range("K64").Value = 121
$ws.Range("K65").Value = 423
$ws.Range("K67").Value = 709
$ws.Range("K69").Value = 41
$ws.Range("K72").Value = 90
$ws.Range("K73").Value = 162
$ws.Range("K77").Value = 130
$ws.Range("K78").Value = 214
$ws.Range("K79").Value = 464
$ws.Range("K80").Value = 66
$ws.Range("K83").Value = 416
$ws.Range("K84").Value = 141
$ws.Range("K85").Value = 880
$ws.Range("K86").Value = 125
$ws.Range("K88").Value = 204
$ws.Range("K90").Value = 167
$ws.Range("K91").Value = 207
$ws.Range("K94").Value = 251
$ws.Range("K95").Value = 315
$ws.Range("K99").Value = 314
$ws.Range("K100").Value = 36
$ws.Range("K101").Value = 18634

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 202
$ws.Range("K3").Value = 255
$ws.Range("K7").Value = 709

# South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 141

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 285
$ws.Range("K3").Value = 359
$ws.Range("K6").Value = 282
$ws.Range("K7").Value = 1001

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 235

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 173
$ws.Range("K7").Value = 547

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 157
$ws.Range("K5").Value = 3

# Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K7").Value = 98

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 181
$ws.Range("K4").Value = 27
$ws.Range("K6").Value = 262
$ws.Range("K7").Value = 690

# Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 104

# Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 214

# Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 208

# Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 194

# Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 207

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 148
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 464

# Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 121

# Calumet Heights
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 126

# Wrigleyville
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 36

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 187
$ws.Range("K3").Value = 177
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 554

# West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 47
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 251

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 127

# Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 52

# Greektown
$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 23

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 118
$ws.Range("K3").Value = 93
$ws.Range("K7").Value = 353

# Gold Coast
$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 29

# Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 75

# Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 53
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 162

# United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 204

# Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 44

# Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 125

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 167

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 165

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 298
$ws.Range("K5").Value = 26
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 880

# Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 90

# Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 130

# Rush & Division
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 66

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 132
$ws.Range("K4").Value = 27
$ws.Range("K7").Value = 488

# Mount Greenwood
$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 16
